# Tekken 8 frame data update:
#  1. Flip the sign of every numeric value in column E ("Block") for rows 2-96.
#  2. Remove the duplicated tail rows 97-131 (leftover from a scraper bug),
#     which shrinks the sheet dimension from A1:H131 down to A1:H96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -is [double]) {
        $cell.Value2 = (0 - $val)
    }
}

$ws.Range("97:131").Delete()
